# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Gilgamesh_Profits workbook (per commit "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""

$ws.Range("H129").Value = 2523.1333
$ws.Range("I129").Value = 932
$ws.Range("K129").Value = 2796
$ws.Range("M129").Value = 2204

$ws.Range("H132").Value = 7720.8335
$ws.Range("I132").Value = 4465.2573
$ws.Range("J132").Value = 23998.715
$ws.Range("K132").Value = 13395.7719
$ws.Range("L132").Value = 71996.145
$ws.Range("M132").Value = -10865.7719
$ws.Range("N132").Value = -77056.145

$ws.Range("H134").Value = 78786.73
$ws.Range("J134").Value = 78786.73
$ws.Range("L134").Value = 78786.73
$ws.Range("N134").Value = -88926.73

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 736.64
$ws.Range("I2").Value = 484.88235
$ws.Range("J2").Value = 1271.625
$ws.Range("K2").Value = 484.88235
$ws.Range("L2").Value = 1271.625
$ws.Range("M2").Value = -371.88235
$ws.Range("N2").Value = -1497.625

$ws.Range("H32").Value = 4734.9067
$ws.Range("I32").Value = 3653.5483
$ws.Range("J32").Value = 7528.4165
$ws.Range("K32").Value = 3653.5483
$ws.Range("L32").Value = 7528.4165
$ws.Range("M32").Value = -3366.5483
$ws.Range("N32").Value = -8102.4165

$ws.Range("H45").Value = 30318.467
$ws.Range("I45").Value = 34057.617
$ws.Range("K45").Value = 34057.617
$ws.Range("M45").Value = -33680.617

$ws.Range("H74").Value = 243424.34
$ws.Range("I74").Value = 309968.34
$ws.Range("J74").Value = 3866
$ws.Range("K74").Value = 309968.34
$ws.Range("L74").Value = 3866
$ws.Range("M74").Value = -309094.34
$ws.Range("N74").Value = -5614

$ws.Range("H77").Value = 243424.34
$ws.Range("I77").Value = 309968.34
$ws.Range("J77").Value = 3866
$ws.Range("K77").Value = 1549841.7
$ws.Range("L77").Value = 19330
$ws.Range("M77").Value = -1545473.7
$ws.Range("N77").Value = -28066

$ws.Range("H88").Value = 3367.8462
$ws.Range("I88").Value = 1784.6
$ws.Range("K88").Value = 1784.6
$ws.Range("M88").Value = -1378.6

$ws.Range("H91").Value = 3367.8462
$ws.Range("I91").Value = 1784.6
$ws.Range("K91").Value = 1784.6
$ws.Range("M91").Value = -380.5999999999999

$ws.Range("H110").Value = 2998.0435
$ws.Range("I110").Value = 1665.3334
$ws.Range("K110").Value = 1665.3334
$ws.Range("M110").Value = 379.6666

$ws.Range("H116").Value = 736.64
$ws.Range("I116").Value = 484.88235
$ws.Range("J116").Value = 1271.625
$ws.Range("K116").Value = 484.88235
$ws.Range("L116").Value = 1271.625
$ws.Range("M116").Value = 1809.11765
$ws.Range("N116").Value = -5859.625

$ws.Range("H132").Value = 1824.0444
$ws.Range("J132").Value = 4080.1667
$ws.Range("L132").Value = 12240.5001
$ws.Range("N132").Value = -17300.5001

$ws.Range("H135").Value = 114996.5
$ws.Range("J135").Value = 114996.5
$ws.Range("L135").Value = 114996.5
$ws.Range("N135").Value = -125136.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 736.64
$ws.Range("I3").Value = 484.88235
$ws.Range("J3").Value = 1271.625
$ws.Range("K3").Value = 484.88235
$ws.Range("L3").Value = 1271.625
$ws.Range("M3").Value = -370.88235
$ws.Range("N3").Value = -1499.625

$ws.Range("H94").Value = 285716830
$ws.Range("J94").Value = 4666.3335
$ws.Range("L94").Value = 4666.3335
$ws.Range("N94").Value = -5568.3335

$ws.Range("H99").Value = 8242.75
$ws.Range("I99").Value = 3884.5715
$ws.Range("K99").Value = 3884.5715
$ws.Range("M99").Value = -2386.5715

$ws.Range("H107").Value = 1557.1578
$ws.Range("I107").Value = 1293.091
$ws.Range("J107").Value = 1920.25
$ws.Range("K107").Value = 1293.091
$ws.Range("L107").Value = 1920.25
$ws.Range("M107").Value = 626.9090000000001
$ws.Range("N107").Value = -5760.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3515.3
$ws.Range("I31").Value = 2156.4358
$ws.Range("J31").Value = 8333.091
$ws.Range("K31").Value = 2156.4358
$ws.Range("L31").Value = 8333.091
$ws.Range("M31").Value = -1861.4358
$ws.Range("N31").Value = -8923.091

$ws.Range("H34").Value = 3515.3
$ws.Range("I34").Value = 2156.4358
$ws.Range("J34").Value = 8333.091
$ws.Range("K34").Value = 2156.4358
$ws.Range("L34").Value = 8333.091
$ws.Range("M34").Value = -1954.4358
$ws.Range("N34").Value = -8737.091

$ws.Range("H132").Value = 2667.1177
$ws.Range("I132").Value = 2260.7273
$ws.Range("J132").Value = 3412.1667
$ws.Range("K132").Value = 6782.1819
$ws.Range("L132").Value = 10236.5001
$ws.Range("M132").Value = -4252.1819
$ws.Range("N132").Value = -15296.5001

$ws.Range("H134").Value = 2746.2693
$ws.Range("J134").Value = 2685.4
$ws.Range("L134").Value = 8056.200000000001
$ws.Range("N134").Value = -13126.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 4333
$ws.Range("J117").Value = 4999.5
$ws.Range("L117").Value = 14998.5
$ws.Range("N117").Value = -21882.5

$ws.Range("H139").Value = 2456.1482
$ws.Range("I139").Value = 795.1539
$ws.Range("K139").Value = 2385.4617
$ws.Range("M139").Value = 2754.5383

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 129543.56
$ws.Range("I70").Value = 170991.5
$ws.Range("K70").Value = 170991.5
$ws.Range("M70").Value = -170721.5

$ws.Range("H73").Value = 129543.56
$ws.Range("I73").Value = 170991.5
$ws.Range("K73").Value = 170991.5
$ws.Range("M73").Value = -170055.5

$ws.Range("H97").Value = 3217.077
$ws.Range("I97").Value = 2504.875
$ws.Range("K97").Value = 2504.875
$ws.Range("M97").Value = -2008.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3437
$ws.Range("I122").Value = 4128.6665
$ws.Range("K122").Value = 12385.9995
$ws.Range("M122").Value = -9935.999500000002

$ws.Range("H132").Value = 6213.0557
$ws.Range("I132").Value = 1372.1666
$ws.Range("K132").Value = 4116.4998
$ws.Range("M132").Value = -1586.4998

$ws.Range("H136").Value = 5454.364
$ws.Range("I136").Value = 6285.5713
$ws.Range("J136").Value = 3999.75
$ws.Range("K136").Value = 18856.7139
$ws.Range("L136").Value = 11999.25
$ws.Range("M136").Value = -16306.7139
$ws.Range("N136").Value = -17099.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 642.1111
$ws.Range("J113").Value = 747.3333
$ws.Range("L113").Value = 2241.9999
$ws.Range("N113").Value = -6581.9999

$ws.Range("H126").Value = 1152.4
$ws.Range("I126").Value = 1152.4
$ws.Range("K126").Value = 3457.2
$ws.Range("M126").Value = -987.2000000000003

$ws.Range("H132").Value = 3321.257
$ws.Range("I132").Value = 3600.138
$ws.Range("J132").Value = 1973.3334
$ws.Range("K132").Value = 10800.414
$ws.Range("L132").Value = 5920.0002
$ws.Range("M132").Value = -8270.414000000001
$ws.Range("N132").Value = -10980.0002

$ws.Range("H133").Value = 94500
$ws.Range("J133").Value = 94500
$ws.Range("L133").Value = 94500
$ws.Range("N133").Value = -104620

$ws.Range("H136").Value = 5108
$ws.Range("I136").Value = 1893.4615
$ws.Range("K136").Value = 5680.3845
$ws.Range("M136").Value = -3130.3845
